$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New group_words values for rows 67-87 (B column)
$values = @(
    "Viện KSND , truy tố , tham nhũng",
    "Viện KSND , truy tố , gây thất thoát , tài sản nhà nước",
    "Viện KSND , truy tố , vi phạm quy định , gây hậu quả nghiêm trọng",
    "Viện KSND , truy tố , lợi dụng chức vụ quyền hạn",
    "Viện KSND , truy tố , hối lộ",
    "Viện KSND , tham ô tài sản, truy tố",
    "Viện KSND , nhận hối lộ, truy tố",
    "Viện KSND , lạm dụng chức vụ quyền hạn, truy tố",
    "Viện KSND , giả mạo trong công tác vì vụ lợi, truy tố",
    "Viện KSND , nhũng nhiễu vì vụ lợi, truy tố",
    "Viện KSND , khởi tố , tham nhũng",
    "Viện KSND , khởi tố , gây thất thoát , tài sản nhà nước",
    "Viện KSND , khởi tố , vi phạm quy định , gây hậu quả nghiêm trọng",
    "Viện KSND , khởi tố , lợi dụng chức vụ quyền hạn",
    "Viện KSND , khởi tố , hối lộ",
    "Viện KSND , tham ô tài sản, khởi tố",
    "Viện KSND , nhận hối lộ, khởi tố",
    "Viện KSND , lạm dụng chức vụ quyền hạn, khởi tố",
    "Viện KSND , giả mạo trong công tác vì vụ lợi, khởi tố",
    "Viện KSND , nhũng nhiễu vì vụ lợi, khởi tố",
    "Viện KSND, truy tố, gây thất thu, ngân sách nhà nước"
)

# Row 67: first new row, plain formula referencing A66 (continues existing series)
$ws.Range("A67").Formula = "=1+A66"
$ws.Range("B67").Value = $values[0]

# Rows 68-87: subsequent rows, each referencing the row above
for ($r = 68; $r -le 87; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=1+A$prev"
    $idx = $r - 67
    $ws.Range("B$r").Value = $values[$idx]
}

# B68:B87 are brand-new cells, so they picked up the plain column style
# instead of the existing "group_words" text style used by B2:B67. Re-apply
# that look (dark slate font color) to match the rest of the column.
$ws.Range("B68:B87").Font.Color = 5057303

# Re-create the (now-cleared) AutoFilter over the old extent, then turn it back
# off: this mirrors the leftover hidden _FilterDatabase name without leaving
# a live autoFilter on the sheet.
$ws.Range("A1:B66").AutoFilter() | Out-Null
$ws.AutoFilterMode = $false
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$B`$66")
$fdb.Visible = $false

# Match the final on-screen selection left behind by the edit.
$ws.Range("B79").Select() | Out-Null
